$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95; this shifts the existing rows 95-113
# down to 96-114 (matching the diff's row-shift pattern) and grows the
# sheet dimension from R113 to R114.
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new weekly record.
$ws.Cells.Item(95, 1).Value = 10
$ws.Cells.Item(95, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(95, 3).Value = "La Araucanía"
$ws.Cells.Item(95, 4).Value = 44476
$ws.Cells.Item(95, 5).Value = 9
$ws.Cells.Item(95, 6).Value = 100114007
$ws.Cells.Item(95, 7).Value = "Jengibre"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 50
$ws.Cells.Item(95, 11).Value = 20000
$ws.Cells.Item(95, 12).Value = 20000
$ws.Cells.Item(95, 13).Value = 20000
$ws.Cells.Item(95, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(95, 15).Value = "Perú"
$ws.Cells.Item(95, 16).Value = 1538
$ws.Cells.Item(95, 17).Value = 13
$ws.Cells.Item(95, 18).Value = "Hortaliza"
